# Insert a new data row at row 46 (pushing existing rows 46-154 down to 47-155),
# then populate the new row 46 with a copy of what is now row 47 (the original
# row 46 data), and finally update its Fecha (column D) to the new date value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("46:46").Insert()
$ws.Range("A47:R47").Copy($ws.Range("A46:R46"))
$ws.Cells.Item(46, 4).Value = 45002
